$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("B1").Value = "toplam"
    $ws.Range("C1").Value = "erkek"
    $ws.Range("D1").Value = "kadın"
    $ws.Range("A2").Value = "toplam"

    for ($row = 3; $row -le 83; $row++) {
        $ws.Cells.Item($row, 1).Value = $row - 2
    }
}
